# Update gh-pages to output generated at 456a3b4
# Applies value updates to 展览(Exhibition), 演出(Performance), 全部类型(All Types)
# sheets and inserts a new row into 本地生活(Local Life) sheet.

$wb = $excel.ActiveWorkbook

# --- Sheet "展览" (Exhibition) ---
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F6").Value = 232
$ws1.Range("F7").Value = 12841
$ws1.Range("F8").Value = 42
$ws1.Range("F9").Value = 115
$ws1.Range("F10").Value = 208
$ws1.Range("F11").Value = 2904
$ws1.Range("F13").Value = 6199
$ws1.Range("F14").Value = 61
$ws1.Range("F16").Value = 3297
$ws1.Range("F19").Value = 120
$ws1.Range("F24").Value = 3544
$ws1.Range("F25").Value = 76
$ws1.Range("F27").Value = 2663
$ws1.Range("F28").Value = 359
$ws1.Range("F29").Value = 1849
$ws1.Range("F31").Value = 192
$ws1.Range("F32").Value = 6450
$ws1.Range("F35").Value = 125
$ws1.Range("F36").Value = 1937
$ws1.Range("F38").Value = 80
$ws1.Range("F39").Value = 1002
$ws1.Range("F41").Value = 197
$ws1.Range("F42").Value = 213
$ws1.Range("F43").Value = 1143
$ws1.Range("F44").Value = 1130
$ws1.Range("F46").Value = 1169
$ws1.Range("F47").Value = 1715
$ws1.Range("F49").Value = 1160

# --- Sheet "演出" (Performance) ---
$ws2 = $wb.Worksheets.Item("演出")
$ws2.Range("F10").Value = 119
$ws2.Range("F14").Value = 924
$ws2.Range("F16").Value = 91

# --- Sheet "本地生活" (Local Life) ---
$ws3 = $wb.Worksheets.Item("本地生活")
$ws3.Range("F2").Value = 407
$ws3.Range("F3").Value = 571

# Insert new row 4 (new event entry), mirroring the formatting of row 3
$ws3.Range("A3").Copy($ws3.Range("A4"))
$ws3.Range("A4").Value = 3
$ws3.Range("B4").NumberFormat = "@"
$ws3.Range("B4").Value = "2024-07-17"
$ws3.Range("B4").ClearFormats()
$ws3.Range("C4").Value = "北京·“狐妖小红娘”限时快闪店"
$ws3.Range("D4").Value = "王府井大街88号 北京王府井银泰in88购物中心"
$ws3.Range("E4").Value = "2024.07.17 10:00-10.31 22:00"
$ws3.Range("F4").Value = 1
$ws3.Range("G4").Value = 98
$ws3.Range("H4").Value = "https://show.bilibili.com/platform/detail.html?id=89613"
$ws3.Range("I4").Value = "//i0.hdslb.com/bfs/openplatform/202407/n3TXriJX1721203778030.jpeg"

# --- Sheet "全部类型" (All Types) ---
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F6").Value = 407
$ws4.Range("F7").Value = 571
$ws4.Range("F8").Value = 232
$ws4.Range("F11").Value = 12841
$ws4.Range("F12").Value = 115
$ws4.Range("F14").Value = 208
$ws4.Range("F15").Value = 2904
$ws4.Range("F17").Value = 6199
$ws4.Range("F18").Value = 61
$ws4.Range("F23").Value = 119
$ws4.Range("F24").Value = 3544
$ws4.Range("F25").Value = 76
$ws4.Range("F27").Value = 2664
$ws4.Range("F28").Value = 1849
$ws4.Range("F30").Value = 192
$ws4.Range("F31").Value = 6450
$ws4.Range("F32").Value = 91
$ws4.Range("F34").Value = 125
$ws4.Range("F35").Value = 1937
$ws4.Range("F38").Value = 80
$ws4.Range("F39").Value = 1002
$ws4.Range("F40").Value = 197
$ws4.Range("F41").Value = 213
$ws4.Range("F42").Value = 1143
$ws4.Range("F43").Value = 1130
$ws4.Range("F45").Value = 1169
$ws4.Range("F47").Value = 1715
$ws4.Range("F49").Value = 1160

"edit applied"
